$wb = $excel.ActiveWorkbook

# Helper: write $value into $range as literal TEXT (preserving leading
# zeros / exact formatting like "008707" or "0.71") without leaving a
# lingering explicit NumberFormat-driven style on the cell - matching
# the plain (un-styled) text cells used throughout this workbook.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# ---------------------------------------------------------------------
# 1) Insert a new "2022-Q1" sheet (fund holdings detail), positioned
#    right before the "总计" (Total) summary sheet. We copy the most
#    recent quarter sheet ("2021-Q4") as a template so styles, number
#    formats and cell types (e.g. text-typed numeric-looking strings)
#    match the existing quarterly sheets exactly, then overwrite the
#    data cells with the 2022-Q1 numbers.
# ---------------------------------------------------------------------
$template = $wb.Worksheets.Item("2021-Q4")
$total = $wb.Worksheets.Item("总计")
$template.Copy($total)

# Note: worksheet handles returned by Item() are position-based, so
# after the Copy() shuffles sheet order we must re-fetch by name
# instead of reusing the old $total handle.
$newSheet = $wb.Worksheets.Item("2021-Q4 (2)")
$newSheet.Name = "2022-Q1"

# Fund holding rows for 2022-Q1.
#   code, fund name, 基金规模, 股票总仓位, 仓位占比, 持有市值(亿元), 仓位排名
$fundRows = @(
    @("539003", "建信富时100指数（QDII）人民币A", "0.71", "92.86", "4.64", "0.0329", 8),
    @("008707", "建信富时100指数（QDII）美元现汇A", "0.71", "92.86", "4.64", "0.0329", 8),
    @("008706", "建信富时100指数（QDII）人民币C", "0.20", "92.86", "4.64", "0.0093", 8),
    @("008708", "建信富时100指数（QDII）美元现汇C", "0.20", "92.86", "4.64", "0.0093", 8)
)

for ($i = 0; $i -lt $fundRows.Count; $i++) {
    $r = $i + 2
    $row = $fundRows[$i]
    $newSheet.Range("A$r").Value = $i
    Set-TextValue $newSheet.Range("B$r") $row[0]
    Set-TextValue $newSheet.Range("C$r") $row[1]
    Set-TextValue $newSheet.Range("D$r") $row[2]
    Set-TextValue $newSheet.Range("E$r") $row[3]
    Set-TextValue $newSheet.Range("F$r") $row[4]
    Set-TextValue $newSheet.Range("G$r") $row[5]
    $newSheet.Range("H$r").Value = $row[6]
}

# ---------------------------------------------------------------------
# 2) Add the 2022-Q1 summary row to the "总计" sheet, above the
#    existing rows (most-recent-first ordering), shifting everything
#    else down by one row and renumbering the index column (A).
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()

# New row 2 loses its style during the insert (it inherits the header
# row's formatting) - clear it and re-apply the "index" style (the one
# used by the other A-column cells) from row 3.
$total.Range("A2:D2").ClearFormats()
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 4
$total.Range("D2").Value = 0.08

# Renumber the index column for the rows that shifted down.
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3
$total.Range("A6").Value = 4
$total.Range("A7").Value = 5
